$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-parsed as a number by Excel (so formatting like trailing zeros survives)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '51.667.07'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '2.956.82'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '378.37'
$ws.Range("E5").Value = '  +7.52%  '
$ws.Range("D6").Value = '104.55'
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '0.544'
$ws.Range("E7").Value = '  -0.63%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.598'
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("D10").Value = '37.27'
$ws.Range("E10").Value = '  +0.16%  '
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").Value = '0.0842'
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").Value = '18.54'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").Value = '3.425.13'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '7.44'
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").Value = '2.957.79'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '0.954'
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").Value = '51.654.86'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("E19").Value = '  +4.91%  '
$ws.Range("D20").Value = '7.42'
$ws.Range("E20").Value = '  +2.20%  '
$ws.Range("D21").Value = '13.15'
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = '0.0₃0956'
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("D23").Value = '68.60'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '262.99'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = '2.79'
$ws.Range("E25").Value = '  +4.69%  '
$ws.Range("D26").Value = '7.40'
$ws.Range("E26").Value = '  +19.38%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '4.16'
$ws.Range("E27").Value = '  -4.00%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '0.169'
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '7.41'
$ws.Range("E30").Value = '  +3.38%  '
$ws.Range("D31").Value = '25.96'
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").Value = '0.103'
$ws.Range("E32").Value = '  -4.13%  '
$ws.Range("D33").Value = '9.89'
$ws.Range("E33").Value = '  -0.95%  '
$ws.Range("D34").Value = '52.58'
$ws.Range("E34").Value = '  +3.63%  '
$ws.Range("D35").Value = '34.30'
$ws.Range("E35").Value = '  -2.47%  '
$ws.Range("E36").Value = '  -4.10%  '
$ws.Range("E37").Value = '  +3.17%  '
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("D39").Value = '3.05'
$ws.Range("E39").Value = '  -4.94%  '
$ws.Range("D40").Value = '2.65'
$ws.Range("E40").Value = '  -4.99%  '
$ws.Range("D41").Value = '17.33'
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").Value = '1.84'
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("D44").Value = '123.69'
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("D45").Value = '21.99'
$ws.Range("E45").Value = '  -2.95%  '
$ws.Range("D46").Value = '0.283'
$ws.Range("E46").Value = '  +20.22%  '
$ws.Range("D47").Value = '2.08'
$ws.Range("E47").Value = '  -3.70%  '
$ws.Range("D48").Value = '2.029.49'
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("D49").Value = '2.32'
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").Value = '3.21'
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("D51").Value = '0.0331'
$ws.Range("E51").Value = '  +3.33%  '
